# "Generate Report for Handback" - refresh the handback/handoff timestamps
# recorded on the Overview, zh-cn and de-de sheets of the handback status
# report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-20 11:08:01"

# zh-cn sheet, first data row: Correspond Handoff Datetime / Correspond
# Handback DateTime.
$wsZhCn.Range("H2").Value = "2016-08-20 11:07:56"
$wsZhCn.Range("K2").Value = "2016-08-20 11:08:26"

# de-de sheet, first data row: Correspond Handoff Datetime (shares the
# same timestamp as the Overview sheet) / Correspond Handback DateTime.
$wsDeDe.Range("H2").Value = "2016-08-20 11:08:01"
$wsDeDe.Range("K2").Value = "2016-08-20 11:08:33"
